$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Format column D as Text before writing numeric-looking price strings,
# then restore the Normal style so no new cell style is introduced.
$colD = $ws.Range("D2:D51")
$colD.NumberFormat = "@"

$ws.Range('D2').Value = '36.757.82'
$ws.Range('E2').Value = '  +3.71%  '
$ws.Range('D3').Value = '1.911.64'
$ws.Range('E3').Value = '  +1.34%  '
$ws.Range('E4').Value = '  +0.01%  '
$ws.Range('D5').Value = '250.06'
$ws.Range('E5').Value = '  +1.31%  '
$ws.Range('D6').Value = '0.698'
$ws.Range('E6').Value = '  +0.02%  '
$ws.Range('E7').Value = '  +0.02%  '
$ws.Range('D8').Value = '46.55'
$ws.Range('E8').Value = '  +6.93%  '
$ws.Range('E9').Value = '  +4.71%  '
$ws.Range('D10').Value = '58.02'
$ws.Range('E10').Value = '  +8.03%  '
$ws.Range('D11').Value = '0.0758'
$ws.Range('E11').Value = '  +1.31%  '
$ws.Range('D12').Value = '0.0999'
$ws.Range('E12').Value = '  +2.04%  '
$ws.Range('D13').Value = '14.60'
$ws.Range('E13').Value = '  +7.35%  '
$ws.Range('D14').Value = '0.814'
$ws.Range('E14').Value = '  +5.24%  '
$ws.Range('D15').Value = '2.191.45'
$ws.Range('E15').Value = '  +1.42%  '
$ws.Range('D16').Value = '5.11'
$ws.Range('E16').Value = '  +2.95%  '
$ws.Range('D17').Value = '1.915.97'
$ws.Range('E17').Value = '  +1.21%  '
$ws.Range('D18').Value = '36.765.41'
$ws.Range('E18').Value = '  +3.71%  '
$ws.Range('D19').Value = '74.44'
$ws.Range('E19').Value = '  +1.13%  '
$ws.Range('D20').Value = '0.0₃0854'
$ws.Range('E20').Value = '  +3.21%  '
$ws.Range('D21').Value = '13.58'
$ws.Range('E21').Value = '  +5.54%  '
$ws.Range('D22').Value = '250.46'
$ws.Range('E22').Value = '  +2.13%  '
$ws.Range('D23').Value = '5.16'
$ws.Range('E23').Value = '  -0.61%  '
$ws.Range('B24').Value = 'Dai'
$ws.Range('C24').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D24').Value = '1.00'
$ws.Range('E24').Value = '  +0.01%  '
$ws.Range('B25').Value = 'Toncoin'
$ws.Range('C25').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D25').Value = '2.51'
$ws.Range('E25').Value = '  -4.64%  '
$ws.Range('D26').Value = '2.23'
$ws.Range('E26').Value = '  +3.36%  '
$ws.Range('D27').Value = '167.35'
$ws.Range('E27').Value = '  +1.49%  '
$ws.Range('D28').Value = '8.77'
$ws.Range('E28').Value = '  +1.21%  '
$ws.Range('D29').Value = '18.67'
$ws.Range('E29').Value = '  +1.62%  '
$ws.Range('E30').Value = '  +0.07%  '
$ws.Range('D31').Value = '4.63'
$ws.Range('E31').Value = '  +7.55%  '
$ws.Range('D32').Value = '0.0616'
$ws.Range('E32').Value = '  +3.06%  '
$ws.Range('D33').Value = '4.32'
$ws.Range('E33').Value = '  +2.97%  '
$ws.Range('B34').Value = 'Kaspa'
$ws.Range('C34').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D34').Value = '0.0893'
$ws.Range('E34').Value = '  +21.02%  '
$ws.Range('B35').Value = 'WEMIXToken'
$ws.Range('C35').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D35').Value = '1.91'
$ws.Range('E35').Value = '  +1.59%  '
$ws.Range('E36').Value = '  +0.12%  '
$ws.Range('D37').Value = '18.86'
$ws.Range('E37').Value = '  +54.36%  '
$ws.Range('D38').Value = '1.48'
$ws.Range('E38').Value = '  -1.37%  '
$ws.Range('D39').Value = '0.876'
$ws.Range('E39').Value = '  +2.54%  '
$ws.Range('D40').Value = '2.01'
$ws.Range('E40').Value = '  +2.08%  '
$ws.Range('D41').Value = '104.96'
$ws.Range('E41').Value = '  +7.50%  '
$ws.Range('D42').Value = '0.0230'
$ws.Range('E42').Value = '  +4.81%  '
$ws.Range('D43').Value = '17.86'
$ws.Range('E43').Value = '  +2.91%  '
$ws.Range('E44').Value = '  +21.45%  '
$ws.Range('E45').Value = '  +1.44%  '
$ws.Range('D46').Value = '1.352.31'
$ws.Range('E46').Value = '  +2.88%  '
$ws.Range('D47').Value = '2.40'
$ws.Range('E47').Value = '  -0.44%  '
$ws.Range('D48').Value = '0.0816'
$ws.Range('E48').Value = '  +0.96%  '
$ws.Range('D49').Value = '2.83'
$ws.Range('E49').Value = '  +3.08%  '
$ws.Range('D50').Value = '6.47'
$ws.Range('E50').Value = '  +1.84%  '
$ws.Range('B51').Value = 'RocketPoolETH'
$ws.Range('C51').Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range('D51').Value = '2.091.82'
$ws.Range('E51').Value = '  +1.26%  '

$colD.Style = "Normal"
